# Add more verbose/specific error messages to the data file used for
# import-error testing, and move the active cell selection from D4 to D5.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B: date errors
$ws.Range("B2:B6").Value = "Error: unprocessable date"

# Column C: numeric errors
$ws.Range("C2:C6").Value = "Error: not a number"

# Column D: per-row special combined values
$ws.Range("D2").Value = "Error: wrong number of arguments for field"
$ws.Range("D3").Value = "33;Error: unprocessable date"
$ws.Range("D4").Value = "Error:not a number;18-03-2022"
$ws.Range("D5").Value = "28;Error: unprocessable date"
$ws.Range("D6").Value = "5;02-03-2023"

# Column E: datetime errors
$ws.Range("E2:E6").Value = "Error: unprocessable datetime"

# Column F: time errors
$ws.Range("F2:F6").Value = "Error: unprocessable time"

# Column G: year-related errors (rows 2,3,5 -> out of bounds; rows 4,6 -> not a year)
$ws.Range("G2").Value = "Error: year out of bounds"
$ws.Range("G3").Value = "Error: year out of bounds"
$ws.Range("G4").Value = "Error: not a year"
$ws.Range("G5").Value = "Error: year out of bounds"
$ws.Range("G6").Value = "Error: not a year"

# Column H: sex option errors
$ws.Range("H2:H6").Value = "Error: non-existent option"

# Column I: race option errors (row 5 has two stacked errors)
$ws.Range("I2").Value = "Error: non-existent option"
$ws.Range("I3").Value = "Error: non-existent option"
$ws.Range("I4").Value = "Error: non-existent option"
$ws.Range("I5").Value = "Error: non-existent option;Error: non-existent option"
$ws.Range("I6").Value = "Error: non-existent option"

# Column J: family history option errors (rows 2 and 5 have three stacked errors)
$ws.Range("J2").Value = "Error: non-existent option;Error: non-existent option;Error: non-existent option"
$ws.Range("J3").Value = "Error: non-existent option"
$ws.Range("J4").Value = "Error: non-existent option"
$ws.Range("J5").Value = "Error: non-existent option;Error: non-existent option;Error: non-existent option"
$ws.Range("J6").Value = "Error: non-existent option"

# Move the active selection from D4 to D5, as captured in the saved view.
$ws.Range("D5").Select() | Out-Null
